# Add generated files from README
#
# Inserts a new "Résumé Automation" personal-project entry (bold title,
# descriptive text, and two hyperlinks: "pandoc" and "Résumé Template")
# right after the "Genealogical Family Tree" entry and before the
# "Web Application for Restaurant Employee Scheduling" entry, inside the
# existing "Personal Projects" section.

function Find-FreshRange($doc, [string]$searchText) {
    # Always search from a brand-new Range so we never operate on a
    # stale Range whose Start/End became invalid after a prior edit.
    $r = $doc.Content
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
    if (-not $ok) {
        throw ("Find-FreshRange: text not found: " + $searchText)
    }
    return $r
}

$d = $word.ActiveDocument
$LF = [char]11   # vertical-tab == Word's soft line break (renders as <w:br/>)

# --- 1. Insert the plain-text skeleton of the new paragraph content ----
# Anchor on the run that begins the next project entry and insert the new
# material immediately before it (i.e. right after "...Work in progress"
# and its existing line break for the Genealogical Family Tree project).
$anchor = Find-FreshRange $d "Web Application for Restaurant Employee Scheduling"
$anchor.Collapse(1)

$newBlock = "Résumé Automation" + " " `
    + "GitHub action to export markdown résumé format to PDF, plain text, and JSON. Utilizes" + " " `
    + "pandoc" + " " `
    + "Allows versioning of custom résumés." + " " `
    + "Résumé Template" + " " `
    + "This document was created with this process" `
    + $LF
$anchor.InsertBefore($newBlock)

# --- 2. Bold the new heading run ("Résumé Automation") -----------------
$heading = Find-FreshRange $d "Résumé Automation"
$heading.Font.Bold = $true

# --- 3. Turn the two label runs into real hyperlinks --------------------
$pandocRun = Find-FreshRange $d "pandoc"
$d.Hyperlinks.Add($pandocRun, "https://pandoc.org/") | Out-Null

$templateRun = Find-FreshRange $d "Résumé Template"
$d.Hyperlinks.Add($templateRun, "https://github.com/mikepqr/resume.md") | Out-Null

Write-Output "Résumé Automation project entry inserted."
